$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "decimale token out" setting (row 2) ---
$ws.Range("D2").Value = 2

# --- Update the existing purchase row (row 4) with a new price / amount / date ---
$ws.Range("A4").Value = 0.254
$ws.Range("B4").Value = 103.01
$ws.Range("C4").Value = 45789

# Row 4 no longer carries the "current row" highlight fill, and its
# alignment reverts to the sheet default (general).
$ws.Range("A4:C4").Interior.ColorIndex = -4142
$ws.Range("A4:C4").HorizontalAlignment = 1

# D4 instead now takes on the pale-green highlight.
$ws.Range("D4").Interior.Color = 7072468

# --- Add a new purchase row (row 5) ---
$ws.Range("A5").Value = 0.0001
$ws.Range("B5").Value = 0.000001
$ws.Range("C5").Value = 45961

# --- Recalculate so the summary formulas in row 2 pick up the new data ---
$excel.Calculate() | Out-Null

# --- Restore the selection to D2 ---
$ws.Range("D2").Select() | Out-Null
